# Append two new "Title and Content" slides to the end of the deck,
# matching the style of the existing slides that use that layout
# (e.g. slide 30, "Smoking Cessation"). We duplicate slide 30 as a
# starting point so the new slides inherit the exact placeholder
# geometry / formatting conventions already used throughout this deck,
# then overwrite the text.

$p = $ppt.ActivePresentation

# --- Slide 31: "Team Members - Physicians" -------------------------------

$srcSlide = $p.Slides.Item($p.Slides.Count)
$dup1 = $srcSlide.Duplicate()
$slide31 = $dup1.Item(1)

$slide31.Shapes.Item(1).TextFrame.TextRange.Text = "Team Members - Physicians"

$body31 = $slide31.Shapes.Item(2)
$tr31 = $body31.TextFrame.TextRange
$tr31.Text = "Primary Care Provider`rGastroenterologist`rMedical Oncologist (chemotherapy)`rRadiation Oncologist (radiation)`rSurgeons`rJonathan Salo MD`rJeffrey Hagen MD`rMichael Roach MD"

# First five lines are plain (no-bullet) role headings; the last three
# are regular bulleted names.
$noBullet31 = @(1, 2, 3, 4, 5)
foreach ($i in $noBullet31) {
    $tr31.Paragraphs($i, 1).ParagraphFormat.Bullet.Visible = 0
}

# --- Slide 32: "Team Members - Support Staff" -----------------------------
# Duplicate the same pristine source slide (not slide31) so the new
# paragraphs inherit a uniform starting format rather than picking up
# slide31's already-edited per-paragraph bullet state. Duplicate() inserts
# immediately after its source, so move the new slide to the end
# afterwards to keep it after slide31.

$dup2 = $srcSlide.Duplicate()
$slide32 = $dup2.Item(1)
$slide32.MoveTo($p.Slides.Count)

$slide32.Shapes.Item(1).TextFrame.TextRange.Text = "Team Members - Support Staff"

$body32 = $slide32.Shapes.Item(2)
$tr32 = $body32.TextFrame.TextRange
$tr32.Text = "Dietitian - Liz Koch`rNurses`rMatthew Carpenter RN`rBrandon Galloway LPN`rNavigator - Laura Swift"

# Lines 1, 2 and 5 are plain (no-bullet); lines 3-4 are regular bulleted names.
$noBullet32 = @(1, 2, 5)
foreach ($i in $noBullet32) {
    $tr32.Paragraphs($i, 1).ParagraphFormat.Bullet.Visible = 0
}
